$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 28715
$ws.Range("J87").Value = 29960
$ws.Range("L87").Value = 29960
$ws.Range("N87").Value = -32456

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 28715
$ws.Range("J90").Value = 29960
$ws.Range("L90").Value = 89880
$ws.Range("N90").Value = -102360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4341.9375
$ws.Range("I113").Value = 3966.25
$ws.Range("J113").Value = 4467.1665
$ws.Range("K113").Value = 3966.25
$ws.Range("L113").Value = 4467.1665
$ws.Range("M113").Value = -712.25
$ws.Range("N113").Value = -10975.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1636.0482
$ws.Range("I132").Value = 1308.75
$ws.Range("J132").Value = 3119.8
$ws.Range("K132").Value = 3926.25
$ws.Range("L132").Value = 9359.400000000001
$ws.Range("M132").Value = -1396.25
$ws.Range("N132").Value = -14419.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 879.0357
$ws.Range("I141").Value = 744.52
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 2233.56
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 2946.44
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1464.69
$ws.Range("I32").Value = 1304.0804
$ws.Range("J32").Value = 2539.5386
$ws.Range("K32").Value = 1304.0804
$ws.Range("L32").Value = 2539.5386
$ws.Range("M32").Value = -1017.0804
$ws.Range("N32").Value = -3113.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3392.5334
$ws.Range("I107").Value = 3258.625
$ws.Range("J107").Value = 3545.5715
$ws.Range("K107").Value = 3258.625
$ws.Range("L107").Value = 3545.5715
$ws.Range("M107").Value = -1338.625
$ws.Range("N107").Value = -7385.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3868.0256
$ws.Range("I134").Value = 3363.577
$ws.Range("J134").Value = 4876.923
$ws.Range("K134").Value = 10090.731
$ws.Range("L134").Value = 14630.769
$ws.Range("M134").Value = -7555.731
$ws.Range("N134").Value = -19700.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2875.3542
$ws.Range("I31").Value = 1950.75
$ws.Range("J31").Value = 5649.1665
$ws.Range("K31").Value = 1950.75
$ws.Range("L31").Value = 5649.1665
$ws.Range("M31").Value = -1655.75
$ws.Range("N31").Value = -6239.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2875.3542
$ws.Range("I34").Value = 1950.75
$ws.Range("J34").Value = 5649.1665
$ws.Range("K34").Value = 1950.75
$ws.Range("L34").Value = 5649.1665
$ws.Range("M34").Value = -1748.75
$ws.Range("N34").Value = -6053.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3408.4443
$ws.Range("I58").Value = 3471.2778
$ws.Range("J58").Value = 3157.111
$ws.Range("K58").Value = 3471.2778
$ws.Range("L58").Value = 3157.111
$ws.Range("M58").Value = -3268.2778
$ws.Range("N58").Value = -3563.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 146082.42
$ws.Range("I99").Value = 251877.75
$ws.Range("J99").Value = 5022
$ws.Range("K99").Value = 251877.75
$ws.Range("L99").Value = 5022
$ws.Range("M99").Value = -250379.75
$ws.Range("N99").Value = -8018

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 146082.42
$ws.Range("I126").Value = 251877.75
$ws.Range("J126").Value = 5022
$ws.Range("K126").Value = 755633.25
$ws.Range("L126").Value = 15066
$ws.Range("M126").Value = -753163.25
$ws.Range("N126").Value = -20006

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3408.4443
$ws.Range("I136").Value = 3471.2778
$ws.Range("J136").Value = 3157.111
$ws.Range("K136").Value = 10413.8334
$ws.Range("L136").Value = 9471.332999999999
$ws.Range("M136").Value = -7863.8334
$ws.Range("N136").Value = -14571.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 390.7857
$ws.Range("I5").Value = 390.7857
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1172.3571
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1060.3571
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 481
$ws.Range("I14").Value = 481
$ws.Range("K14").Value = 1443
$ws.Range("M14").Value = -1270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 700
$ws.Range("J62").Value = 3100
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 9300
$ws.Range("M62").Value = -1414
$ws.Range("N62").Value = -10672

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 2700
$ws.Range("I65").Value = 700
$ws.Range("J65").Value = 3100
$ws.Range("K65").Value = 6300
$ws.Range("L65").Value = 27900
$ws.Range("M65").Value = -2868
$ws.Range("N65").Value = -34764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 450
$ws.Range("J68").Value = 450
$ws.Range("L68").Value = 1350
$ws.Range("N68").Value = -2972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 450
$ws.Range("J71").Value = 450
$ws.Range("L71").Value = 4050
$ws.Range("N71").Value = -12162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 390.7857
$ws.Range("I135").Value = 390.7857
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3517.0713
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -982.0713000000001
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1251.3
$ws.Range("I122").Value = 1126.625
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3379.875
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -929.875
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1999.4166
$ws.Range("I126").Value = 1674.8889
$ws.Range("J126").Value = 2973
$ws.Range("K126").Value = 5024.6667
$ws.Range("L126").Value = 8919
$ws.Range("M126").Value = -2554.6667
$ws.Range("N126").Value = -13859

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4313.6665
$ws.Range("I132").Value = 6342.6665
$ws.Range("J132").Value = 3154.238
$ws.Range("K132").Value = 19027.9995
$ws.Range("L132").Value = 9462.714
$ws.Range("M132").Value = -16497.9995
$ws.Range("N132").Value = -14522.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3804
$ws.Range("I7").Value = 3995
$ws.Range("J7").Value = 3740.3333
$ws.Range("K7").Value = 3995
$ws.Range("L7").Value = 3740.3333
$ws.Range("M7").Value = -3883
$ws.Range("N7").Value = -3964.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3218.8333
$ws.Range("I40").Value = 2803.1428
$ws.Range("J40").Value = 3800.8
$ws.Range("K40").Value = 2803.1428
$ws.Range("L40").Value = 3800.8
$ws.Range("M40").Value = -2667.1428
$ws.Range("N40").Value = -4072.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1380.9714
$ws.Range("I93").Value = 1527.2
$ws.Range("J93").Value = 1186
$ws.Range("K93").Value = 1527.2
$ws.Range("L93").Value = 1186
$ws.Range("M93").Value = -279.2
$ws.Range("N93").Value = -3682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3804
$ws.Range("I126").Value = 3995
$ws.Range("J126").Value = 3740.3333
$ws.Range("K126").Value = 11985
$ws.Range("L126").Value = 11220.9999
$ws.Range("M126").Value = -9515
$ws.Range("N126").Value = -16160.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11196.88
$ws.Range("I132").Value = 3661.4443
$ws.Range("J132").Value = 15435.5625
$ws.Range("K132").Value = 10984.3329
$ws.Range("L132").Value = 46306.6875
$ws.Range("M132").Value = -8454.332900000001
$ws.Range("N132").Value = -51366.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2763.4285
$ws.Range("I136").Value = 1458.0536
$ws.Range("J136").Value = 7984.9287
$ws.Range("K136").Value = 4374.1608
$ws.Range("L136").Value = 23954.7861
$ws.Range("M136").Value = -1824.1608
$ws.Range("N136").Value = -29054.7861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 47315
$ws.Range("J140").Value = 47315
$ws.Range("L140").Value = 47315
$ws.Range("N140").Value = -57675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2514
$ws.Range("I126").Value = 2764.3
$ws.Range("J126").Value = 2286.4546
$ws.Range("K126").Value = 8292.900000000001
$ws.Range("L126").Value = 6859.3638
$ws.Range("M126").Value = -5822.900000000001
$ws.Range("N126").Value = -11799.3638
